$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1 = 14, Q1 = 15, matching the formatting of the existing header cells ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# --- Data rows 2-25: update I/K/M/O columns, and add new P/Q columns ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I column: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K column: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M column: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O column: 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P column: new = 2
    $ws.Cells.Item($r, 17).Value = 2   # Q column: new = 2
}
